$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New cell values for row 2 (demodulation technique details per TNC) ---
$ws.Range("B2").Value = "XR-2211 PLL Chip"
$ws.Range("C2").Value = "T3 uses 8-bit samples. The current sample is multiplied with a sample from a fixed delay (I think it should be 1/4 of the 1700 hz AFSK center frequency) and then the result is low-pass filtered. A separate section handles clock recovery."
$ws.Range("D2").Value = "In the OT1+ the samples were quantized to a single bit, but the T3 uses 8-bit samples. The current sample is multiplied with a sample from a fixed delay (I think it should be 1/4 of the 1700 hz AFSK center frequency) and then the result is low-pass filtered. A separate section handles clock recovery."
$ws.Range("E2").Value = "T3 uses 8-bit samples. The current sample is multiplied with a sample from a fixed delay (I think it should be 1/4 of the 1700 hz AFSK center frequency) and then the result is low-pass filtered. A separate section handles clock recovery."
$ws.Range("F2").Value = "The TCM3105 chip demodulator is an edge-triggered multivibrator that triggers off positive- and negative-going edges"
$ws.Range("H2").Value = "PLL XR-2211"
$ws.Range("I2").Value = "Internal AMD 7910 chip uses digital filters"
$ws.Range("J2").Value = "Analog Filters"
$ws.Range("K2").Value = "Analog Filters"
$ws.Range("L2").Value = "Digital Filters"

# --- Row 3 addition ---
$ws.Range("B3").Value = "10mv-3v"

# --- Row 4 (new row): supporting document links ---
$ws.Range("B4").Value = "http://html.alldatasheet.com/html-pdf/80494/EXAR/XR-2211/81/1/XR-2211.html"
$ws.Range("H4").Value = "http://www.radiomanual.info/schemi/ACC_packet/MFJ-1278_MFJ-1278T_user.pdf"
$ws.Range("I4").Value = "http://pdf1.alldatasheet.com/datasheet-pdf/view/124526/AMD/AM7910PC.html"

# --- Row 5 (new row): supporting document link ---
$ws.Range("I5").Value = "http://www.repeater-builder.com/aea/pdf/aea-pk-88-user-manual.pdf"

# G2 is styled (wrapped) but stays empty -- give it the same formatting as the rest of row 2.
$ws.Range("B2:L2").WrapText = $true

# Row 2 grows tall enough to show the wrapped paragraphs.
$ws.Rows.Item(2).RowHeight = 231.75

# Widen columns B:L so the new text is readable.
$ws.Range("B1:L1").ColumnWidth = 20.7

# Selection moves to A3 (matches the last active cell when the edit was saved).
$ws.Range("A3").Select()

# Portrait page orientation for printing the (now much taller) sheet.
$ws.PageSetup.Orientation = 1
